$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: "Dec 25 2018" -> "Jan 2 2019" and move the "_GoBack" bookmark
#          from after "semester" to right after the new date text.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Dec 25 2018", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Jan 2 2019", 2)

$p6 = $d.Paragraphs.Item(6)
$endPos = $p6.Range.End - 1   # position right after "Jan 2 2019", before the paragraph mark

# Work around a zero-width-range placement quirk: temporarily insert a
# placeholder character so the bookmark position is no longer exactly at
# the end-of-paragraph boundary, add the bookmark, then remove the
# placeholder again.
$tempIns = $d.Range($endPos, $endPos)
$tempIns.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$tempDel = $d.Range($endPos, $endPos + 1)
$tempDel.Delete()

# ------------------------------------------------------------------
# Edit 2: "'s kennel" -> "'s Kennel" (capitalise the K)
# ------------------------------------------------------------------
$d.Content.Find.Execute("’s kennel", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "’s Kennel", 2)
